$d = $word.ActiveDocument
$d.Content.LanguageID = 1033
